# Adds one additional day of violent-crime data (2023-03-14) to the
# '2023' (column J) totals across the Citywide Totals sheet, the By
# Neighborhood summary sheet, and every individual neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1253
$ws.Range('J3').Value = 1336
$ws.Range('J4').Value = 294
$ws.Range('J5').Value = 98
$ws.Range('J6').Value = 1772
$ws.Range('J7').Value = 4753

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J2').Value = 21
$ws.Range('J3').Value = 18
$ws.Range('J7').Value = 59

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J2').Value = 6
$ws.Range('J7').Value = 17

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J5').Value = 6
$ws.Range('J7').Value = 162

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 72
$ws.Range('J6').Value = 57
$ws.Range('J7').Value = 175

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 35

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 16
$ws.Range('J4').Value = 2
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 31
$ws.Range('J3').Value = 27
$ws.Range('J7').Value = 117

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J6').Value = 45
$ws.Range('J7').Value = 132
$ws.Range('J8').Value = 295
$ws.Range('J10').Value = 30
$ws.Range('J22').Value = 8
$ws.Range('J25').Value = 26
$ws.Range('J29').Value = 271
$ws.Range('J30').Value = 17
$ws.Range('J31').Value = 35
$ws.Range('J33').Value = 198
$ws.Range('J36').Value = 70
$ws.Range('J37').Value = 162
$ws.Range('J41').Value = 29
$ws.Range('J42').Value = 191
$ws.Range('J47').Value = 42
$ws.Range('J48').Value = 34
$ws.Range('J51').Value = 62
$ws.Range('J52').Value = 105
$ws.Range('J54').Value = 95
$ws.Range('J63').Value = 16
$ws.Range('J64').Value = 28
$ws.Range('J65').Value = 117
$ws.Range('J67').Value = 175
$ws.Range('J73').Value = 43
$ws.Range('J77').Value = 37
$ws.Range('J79').Value = 141
$ws.Range('J80').Value = 12
$ws.Range('J83').Value = 117
$ws.Range('J84').Value = 49
$ws.Range('J85').Value = 207
$ws.Range('J88').Value = 39
$ws.Range('J89').Value = 54
$ws.Range('J90').Value = 54
$ws.Range('J96').Value = 59
$ws.Range('J101').Value = 4753

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J6').Value = 40
$ws.Range('J7').Value = 117

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 48
$ws.Range('J3').Value = 56
$ws.Range('J6').Value = 80
$ws.Range('J7').Value = 198

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 47
$ws.Range('J7').Value = 95

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J3').Value = 102
$ws.Range('J6').Value = 72
$ws.Range('J7').Value = 271

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 34

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 53
$ws.Range('J3').Value = 76
$ws.Range('J6').Value = 58
$ws.Range('J7').Value = 207

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 15
$ws.Range('J7').Value = 45

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J6').Value = 14
$ws.Range('J7').Value = 29

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 33
$ws.Range('J7').Value = 191

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 30

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 50
$ws.Range('J6').Value = 40
$ws.Range('J7').Value = 141

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J6').Value = 14
$ws.Range('J7').Value = 28

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J3').Value = 14
$ws.Range('J6').Value = 33
$ws.Range('J7').Value = 70

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 26
$ws.Range('J6').Value = 41
$ws.Range('J7').Value = 105

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 11
$ws.Range('J7').Value = 26

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J5').Value = 1
$ws.Range('J7').Value = 42

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 43

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 39

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 98
$ws.Range('J3').Value = 99
$ws.Range('J7').Value = 295

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 14
$ws.Range('J6').Value = 24
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 19
$ws.Range('J7').Value = 62

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J6').Value = 2
$ws.Range('J7').Value = 8

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 9
$ws.Range('J7').Value = 37

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('J2').Value = 2
$ws.Range('J7').Value = 12

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 45
$ws.Range('J7').Value = 132
